# Commit: "debugging examples XLSX files for Microsoft Excel"
#
# This script:
#   1. Removes the 4 legacy cell comments (A2:D2) on the "!!Transaction"
#      worksheet, which also drops the now-unused comments part, its
#      VML legacy-drawing reference, and shrinks the sheet's used
#      dimension from A1:D6 to A1:C6.
#   2. Bumps the embedded ObjTables schema/version marker strings from
#      objTablesVersion='1.0.0' to objTablesVersion='1.0.1' in the
#      header cells of all three worksheets.

$wb = $excel.ActiveWorkbook

$wsToc  = $wb.Worksheets.Item(1)   # "!!_Table of contents"
$wsSch  = $wb.Worksheets.Item(2)   # "!!_Schema"
$wsData = $wb.Worksheets.Item(3)   # "!!Transaction"

# --- 1. Delete the data-entry comments on the Transaction sheet ---------
$wsData.Range("A2").Comment.Delete()
$wsData.Range("B2").Comment.Delete()
$wsData.Range("C2").Comment.Delete()
$wsData.Range("D2").Comment.Delete()

# --- 2. Bump the ObjTables version marker embedded in the header text ---
# The cells carrying this text are protected/locked, so each sheet must
# briefly be unprotected to perform the replace, then restored.

$wsToc.Unprotect()
[void]$wsToc.Cells.Replace("objTablesVersion='1.0.0'", "objTablesVersion='1.0.1'")
$wsToc.Protect($null, $true, $true, $true)

$wsSch.Unprotect()
[void]$wsSch.Cells.Replace("objTablesVersion='1.0.0'", "objTablesVersion='1.0.1'")
$wsSch.Protect($null, $true, $true, $true)

$wsData.Unprotect()
[void]$wsData.Cells.Replace("objTablesVersion='1.0.0'", "objTablesVersion='1.0.1'")
$wsData.Protect($null, $true, $true, $true)
